$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table
$cell = $tbl.Cell(10, 4)
$cell.Shape.TextFrame.TextRange.Text = "See checklist first workday in step 4.2"
